# Append the new Argent (Solar Prices) data row for 2025-05-09.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 69

# The source data is stored as literal text (inline strings) for every
# column, exactly like the preceding rows. Pre-format the target cells as
# Text so Excel's auto-detection doesn't silently turn numeric-looking
# values (dates, plain numbers, "5,370"-style thousands) into real
# numbers/dates.
$targetRange = $ws.Range("A" + $row + ":J" + $row)
$targetRange.NumberFormat = "@"

$ws.Range("A" + $row).Value = "2025-05-09"
$ws.Range("B" + $row).Value = "38"
$ws.Range("C" + $row).Value = "37.28"
$ws.Range("D" + $row).Value = "1"
$ws.Range("E" + $row).Value = "0.265"
$ws.Range("F" + $row).Value = "0.09"
$ws.Range("G" + $row).Value = "5,370"
$ws.Range("H" + $row).Value = "8,039"
$ws.Range("I" + $row).Value = "8,089"
$ws.Range("J" + $row).Value = "7.2577"
